$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: India (sst index unchanged) ---
$ws.Range("B6").Value = 1910681
$ws.Range("C6").Value = 4068
$ws.Range("D6").Value = 1282917
$ws.Range("E6").Value = 587908
$ws.Range("G6").Value = 36
$ws.Range("H6").Value = 39856

# --- Row 17: Pakistan (sst index unchanged) ---
$ws.Range("B17").Value = 281136
$ws.Range("C17").Value = 675
$ws.Range("D17").Value = 254286
$ws.Range("E17").Value = 20836
$ws.Range("G17").Value = 15
$ws.Range("H17").Value = 6014

# --- Row 36: Israel (sst index unchanged) ---
$ws.Range("B36").Value = 76642
$ws.Range("C36").Value = 444
$ws.Range("D36").Value = 51329
$ws.Range("E36").Value = 24752

# --- Rows 55 & 56: Kirguistan / Ghana swap places in ranking ---
# Row 55 now shows Kirguistan with fresh data; row 56 now shows Ghana
# with what used to be row 55's data.
$ws.Range("A55").Value = "Kirguistan"
$ws.Range("B55").Value = 38110
$ws.Range("C55").Value = 569
$ws.Range("D55").Value = 29513
$ws.Range("E55").Value = 7159
$ws.Range("G55").Value = 11
$ws.Range("H55").Value = 1438

$ws.Range("A56").Value = "Ghana"
$ws.Range("B56").Value = 37812
$ws.Range("C56").Value = 0
$ws.Range("D56").Value = 34313
$ws.Range("E56").Value = 3308
$ws.Range("G56").Value = 0
$ws.Range("H56").Value = 191

# --- Row 57: Afganistan (sst index unchanged) ---
$ws.Range("B57").Value = 36829
$ws.Range("C57").Value = 47
$ws.Range("D57").Value = 25742
$ws.Range("E57").Value = 9793
$ws.Range("G57").Value = 6
$ws.Range("H57").Value = 1294

# --- Rows 61 & 62: Uzbekistan / Marruecos swap places in ranking ---
# Row 61 now shows Uzbekistan with fresh data; row 62 now shows Marruecos
# with what used to be row 61's data.
$ws.Range("A61").Value = "Uzbekistan"
$ws.Range("B61").Value = 27314
$ws.Range("C61").Value = 267
$ws.Range("D61").Value = 18051
$ws.Range("E61").Value = 9098
$ws.Range("H61").Value = 165

$ws.Range("A62").Value = "Marruecos"
$ws.Range("B62").Value = 27217
$ws.Range("C62").Value = 0
$ws.Range("D62").Value = 19629
$ws.Range("E62").Value = 7171
$ws.Range("H62").Value = 417

# --- Footer timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 5 de Agosto de 2020 a las 08:05"
